$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "APLICACIÓN DE GRID SEARCH" section added to the right of the
#     existing ANN table (column L), documenting the grid-search results
#     for the ANN model (commit: "preprocesado producción memoria, ann
#     pruebas").
#
# Shared-string insertion order matters (it is derived from first-write
# order), so write L6 before L5 to reproduce the target sharedStrings.xml
# ordering (index 13 = "Best Accuracy..." / index 14 = "APLICACIÓN...").
$ws.Range("L6").Value = "Best Accuracy: 98.20 % Best Parameters: {'activation': 'tanh', 'hidden_layer_sizes': (7, 7, 1), 'learning_rate': 'constant', 'solver': 'adam'}"

# Nudge the cell so it gets its own distinct (but visually identical)
# cell-format entry, the same way the source workbook carries a dedicated
# xf for this cell instead of reusing the default one.
$ws.Range("L6").WrapText = $false

$ws.Range("L5").Value = "APLICACIÓN DE GRID SEARCH"

# Match the author's final selection/cursor position.
$ws.Range("I7").Select()
